$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $value) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws 'D2' '56.438.06'
Set-TextCell $ws 'E2' '  -1.24%  '

Set-TextCell $ws 'D3' '2.981.91'
Set-TextCell $ws 'E3' '  -2.63%  '

Set-TextCell $ws 'D4' '0.999'
Set-TextCell $ws 'E4' '  -0.30%  '

Set-TextCell $ws 'D5' '501.42'
Set-TextCell $ws 'E5' '  -2.00%  '

Set-TextCell $ws 'D6' '133.58'
Set-TextCell $ws 'E6' '  +4.80%  '

Set-TextCell $ws 'E7' '  -0.28%  '

Set-TextCell $ws 'D8' '0.428'
Set-TextCell $ws 'E8' '  -0.51%  '

Set-TextCell $ws 'D9' '7.30'
Set-TextCell $ws 'E9' '  +3.33%  '

Set-TextCell $ws 'D10' '0.106'
Set-TextCell $ws 'E10' '  +2.12%  '

Set-TextCell $ws 'D11' '0.351'
Set-TextCell $ws 'E11' '  -2.22%  '

Set-TextCell $ws 'E12' '  +0.23%  '

Set-TextCell $ws 'D13' '3.487.46'
Set-TextCell $ws 'E13' '  -3.70%  '

Set-TextCell $ws 'E14' '  +3.75%  '

Set-TextCell $ws 'D15' '56.372.42'
Set-TextCell $ws 'E15' '  +3.11%  '

Set-TextCell $ws 'D16' '0.0000149'
Set-TextCell $ws 'E16' '  +3.87%  '

Set-TextCell $ws 'D17' '2.974.62'
Set-TextCell $ws 'E17' '  -3.72%  '

Set-TextCell $ws 'D18' '5.69'
Set-TextCell $ws 'E18' '  +3.50%  '

Set-TextCell $ws 'D19' '12.36'
Set-TextCell $ws 'E19' '  -0.38%  '

Set-TextCell $ws 'D20' '7.75'
Set-TextCell $ws 'E20' '  +2.63%  '

Set-TextCell $ws 'D21' '325.35'
Set-TextCell $ws 'E21' '  -0.64%  '

Set-TextCell $ws 'D22' '1.00'
Set-TextCell $ws 'E22' '  +0.12%  '

Set-TextCell $ws 'D23' '0.471'
Set-TextCell $ws 'E23' '  -3.54%  '

Set-TextCell $ws 'D24' '62.12'
Set-TextCell $ws 'E24' '  -5.10%  '

Set-TextCell $ws 'D25' '0.997'
Set-TextCell $ws 'E25' '  -0.30%  '

Set-TextCell $ws 'D26' '0.164'
Set-TextCell $ws 'E26' '  -1.22%  '

Set-TextCell $ws 'D27' '0.0₃0886'
Set-TextCell $ws 'E27' '  +1.67%  '

Set-TextCell $ws 'E28' '  -0.10%  '

Set-TextCell $ws 'D29' '6.47'
Set-TextCell $ws 'E29' '  -0.54%  '

Set-TextCell $ws 'D30' '6.83'
Set-TextCell $ws 'E30' '  +3.83%  '

Set-TextCell $ws 'D31' '1.74'
Set-TextCell $ws 'E31' '  -2.82%  '

Set-TextCell $ws 'D32' '1.17'
Set-TextCell $ws 'E32' '  -4.02%  '

Set-TextCell $ws 'D33' '20.25'
Set-TextCell $ws 'E33' '  -1.73%  '

Set-TextCell $ws 'D34' '155.13'
Set-TextCell $ws 'E34' '  -0.64%  '

Set-TextCell $ws 'D35' '4.44'
Set-TextCell $ws 'E35' '  -2.60%  '

Set-TextCell $ws 'D36' '1.28'
Set-TextCell $ws 'E36' '  -1.98%  '

Set-TextCell $ws 'D37' '5.55'
Set-TextCell $ws 'E37' '  -5.50%  '

Set-TextCell $ws 'D38' '0.0673'
Set-TextCell $ws 'E38' '  +2.87%  '

Set-TextCell $ws 'D39' '22.96'
Set-TextCell $ws 'E39' '  +0.44%  '

Set-TextCell $ws 'D40' '3.013.95'
Set-TextCell $ws 'E40' '  -3.16%  '

Set-TextCell $ws 'D41' '0.999'
Set-TextCell $ws 'E41' '  -0.24%  '

Set-TextCell $ws 'D42' '35.98'
Set-TextCell $ws 'E42' '  -0.58%  '

Set-TextCell $ws 'D43' '0.641'
Set-TextCell $ws 'E43' '  -3.36%  '

Set-TextCell $ws 'D44' '2.234.23'
Set-TextCell $ws 'E44' '  +0.69%  '

Set-TextCell $ws 'D45' '0.989'
Set-TextCell $ws 'E45' '  -4.37%  '

Set-TextCell $ws 'D46' '1.39'
Set-TextCell $ws 'E46' '  +1.67%  '

Set-TextCell $ws 'D47' '3.56'
Set-TextCell $ws 'E47' '  -3.98%  '

Set-TextCell $ws 'B48' 'dogwifhat'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell $ws 'D48' '1.91'
Set-TextCell $ws 'E48' '  +11.98%  '

Set-TextCell $ws 'B49' 'VeChain'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D49' '0.0234'
Set-TextCell $ws 'E49' '  +4.06%  '

Set-TextCell $ws 'D50' '5.75'
Set-TextCell $ws 'E50' '  -3.03%  '

Set-TextCell $ws 'D51' '18.92'
Set-TextCell $ws 'E51' '  -2.98%  '
